$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.135.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.09%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.238.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.31%  "

# Row 4
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "261.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.33%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "84.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +16.13%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.96%  "

# Row 8
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.612"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.86%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +13.58%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0932"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.17%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.27%  "

# Row 13
$ws.Range("E13").Value = "  +3.04%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.568.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.05%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.64%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.232.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.69%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.790"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.99%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.026.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.19%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000105"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.35%  "

# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.85%  "

# Row 21
$ws.Range("B21").Value = "Litecoin"
$ws.Range("C21").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.70%  "

# Row 22
$ws.Range("E22").Value = "  +12.16%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.28%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.58%  "

# Row 25
$ws.Range("E25").Value = "  +0.10%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "40.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +12.48%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.08%  "

# Row 29
$ws.Range("E29").Value = "  +3.71%  "

# Row 30
$ws.Range("E30").Value = "  +0.13%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.22%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0897"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +12.75%  "

# Row 33
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.76%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.99%  "

# Row 35
$ws.Range("E35").Value = "  +9.16%  "

# Row 36
$ws.Range("E36").Value = "  +2.94%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0368"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +13.73%  "

# Row 38
$ws.Range("E38").Value = "  +8.53%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +15.04%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +27.13%  "

# Row 41
$ws.Range("E41").Value = "  +5.22%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.92"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.11%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.53%  "

# Row 44
$ws.Range("E44").Value = "  +4.50%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.72%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.20%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0989"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.60%  "

# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.84%  "

# Row 49
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +28.90%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.450"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.63%  "

# Row 51
$ws.Range("E51").Value = "  +4.60%  "
